$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2:7) contents only, keep formatting untouched.
$ws.Range("A2:T7").ClearContents()

# Rewrite header row (row 1) exactly as before, column by column,
# so the shared-string pool keeps the same header order.
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# Introduce the data (text) strings in the exact order needed so that newly-created
# shared-string entries land in the desired sequence: ECs, FAPs, sCs, Penk, Ogfr.
# Column A (Sending cluster) first...
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"
$ws.Range("A10").Value = "sCs"
# ...then column D (Target cluster)...
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("D10").Value = "sCs"
# ...then column B (Ligand symbol) and column C (Receptor symbol).
$ws.Range("B2").Value = "Penk"
$ws.Range("B3").Value = "Penk"
$ws.Range("B4").Value = "Penk"
$ws.Range("B5").Value = "Penk"
$ws.Range("B6").Value = "Penk"
$ws.Range("B7").Value = "Penk"
$ws.Range("B8").Value = "Penk"
$ws.Range("B9").Value = "Penk"
$ws.Range("B10").Value = "Penk"
$ws.Range("C2").Value = "Ogfr"
$ws.Range("C3").Value = "Ogfr"
$ws.Range("C4").Value = "Ogfr"
$ws.Range("C5").Value = "Ogfr"
$ws.Range("C6").Value = "Ogfr"
$ws.Range("C7").Value = "Ogfr"
$ws.Range("C8").Value = "Ogfr"
$ws.Range("C9").Value = "Ogfr"
$ws.Range("C10").Value = "Ogfr"

# Fill in the numeric columns (E..T) for every row.
# Row 2: ECs -> ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.860262
$ws.Range("H2").Value = 2.580786
$ws.Range("I2").Value = 0.007129397360689299
$ws.Range("J2").Value = 0.007129397360689299
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.736318
$ws.Range("N2").Value = 32.208954
$ws.Range("O2").Value = 0.5380467806526527
$ws.Range("P2").Value = 0.5380467806526528
$ws.Range("Q2").Value = 9.236046395315999
$ws.Range("R2").Value = 83.12441755784398
$ws.Range("S2").Value = 0.003835949297912396
$ws.Range("T2").Value = 0.003835949297912397
# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.860262
$ws.Range("H3").Value = 2.580786
$ws.Range("I3").Value = 0.007129397360689299
$ws.Range("J3").Value = 0.007129397360689299
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.491314333333333
$ws.Range("N3").Value = 13.473943
$ws.Range("O3").Value = 0.2250806298722816
$ws.Range("P3").Value = 0.2250806298722816
$ws.Range("Q3").Value = 3.863707051021999
$ws.Range("R3").Value = 34.77336345919799
$ws.Range("S3").Value = 0.00160468924855373
$ws.Range("T3").Value = 0.00160468924855373
# Row 4: ECs -> sCs
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.860262
$ws.Range("H4").Value = 2.580786
$ws.Range("I4").Value = 0.007129397360689299
$ws.Range("J4").Value = 0.007129397360689299
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.726614000000001
$ws.Range("N4").Value = 14.179842
$ws.Range("O4").Value = 0.2368725894750656
$ws.Range("P4").Value = 0.2368725894750656
$ws.Range("Q4").Value = 4.066126412868
$ws.Range("R4").Value = 36.595137715812
$ws.Range("S4").Value = 0.001688758814223173
$ws.Range("T4").Value = 0.001688758814223173
# Row 5: FAPs -> ECs
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 118.8986256666667
$ws.Range("H5").Value = 356.695877
$ws.Range("I5").Value = 0.9853690480545674
$ws.Range("J5").Value = 0.9853690480545675
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.736318
$ws.Range("N5").Value = 32.208954
$ws.Range("O5").Value = 0.5380467806526527
$ws.Range("P5").Value = 0.5380467806526528
$ws.Range("Q5").Value = 1276.533454920295
$ws.Range("R5").Value = 11488.80109428266
$ws.Range("S5").Value = 0.5301746440605289
$ws.Range("T5").Value = 0.5301746440605292
# Row 6: FAPs -> FAPs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 118.8986256666667
$ws.Range("H6").Value = 356.695877
$ws.Range("I6").Value = 0.9853690480545674
$ws.Range("J6").Value = 0.9853690480545675
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.491314333333333
$ws.Range("N6").Value = 13.473943
$ws.Range("O6").Value = 0.2250806298722816
$ws.Range("P6").Value = 0.2250806298722816
$ws.Range("Q6").Value = 534.0111016703345
$ws.Range("R6").Value = 4806.099915033011
$ws.Range("S6").Value = 0.2217874859927726
$ws.Range("T6").Value = 0.2217874859927726
# Row 7: FAPs -> sCs
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 118.8986256666667
$ws.Range("H7").Value = 356.695877
$ws.Range("I7").Value = 0.9853690480545674
$ws.Range("J7").Value = 0.9853690480545675
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.726614000000001
$ws.Range("N7").Value = 14.179842
$ws.Range("O7").Value = 0.2368725894750656
$ws.Range("P7").Value = 0.2368725894750656
$ws.Range("Q7").Value = 561.987908656826
$ws.Range("R7").Value = 5057.891177911434
$ws.Range("S7").Value = 0.2334069180012658
$ws.Range("T7").Value = 0.2334069180012658
# Row 8: sCs -> ECs
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.9051680000000001
$ws.Range("H8").Value = 2.715504
$ws.Range("I8").Value = 0.007501554584743267
$ws.Range("J8").Value = 0.007501554584743267
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.736318
$ws.Range("N8").Value = 32.208954
$ws.Range("O8").Value = 0.5380467806526527
$ws.Range("P8").Value = 0.5380467806526528
$ws.Range("Q8").Value = 9.718171491424
$ws.Range("R8").Value = 87.463543422816
$ws.Range("S8").Value = 0.004036187294211261
$ws.Range("T8").Value = 0.004036187294211262
# Row 9: sCs -> FAPs
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.9051680000000001
$ws.Range("H9").Value = 2.715504
$ws.Range("I9").Value = 0.007501554584743267
$ws.Range("J9").Value = 0.007501554584743267
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.491314333333333
$ws.Range("N9").Value = 13.473943
$ws.Range("O9").Value = 0.2250806298722816
$ws.Range("P9").Value = 0.2250806298722816
$ws.Range("Q9").Value = 4.065394012474666
$ws.Range("R9").Value = 36.588546112272
$ws.Range("S9").Value = 0.001688454630955316
$ws.Range("T9").Value = 0.001688454630955316
# Row 10: sCs -> sCs
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9051680000000001
$ws.Range("H10").Value = 2.715504
$ws.Range("I10").Value = 0.007501554584743267
$ws.Range("J10").Value = 0.007501554584743267
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.726614000000001
$ws.Range("N10").Value = 14.179842
$ws.Range("O10").Value = 0.2368725894750656
$ws.Range("P10").Value = 0.2368725894750656
$ws.Range("Q10").Value = 4.278379741152001
$ws.Range("R10").Value = 38.505417670368
$ws.Range("S10").Value = 0.001776912659576688
$ws.Range("T10").Value = 0.001776912659576688
